# Regenerate merged AHB files
#
# The "Gruppe" block spanning rows 62-85 had already been processed
# elsewhere in the sheet (see rows 43/47/51/55 for the pattern): the first
# ("header") row of each Datenelement group gets a gray fill across the
# whole row (with the 2nd column - the field name - in bold), and the
# "Änderung" column (L) loses its bold/gold "ÄNDERUNG" flag, becoming a
# plain empty gray, centered cell - same as every other already-merged
# group. This mirrors that same transformation onto rows 62-85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Donor rows that already carry the target formatting, used as a
# "format painter" source via Copy + PasteSpecial(xlPasteFormats) so the
# exact existing style entries (not new near-duplicate ones) get reused.
$donorHeaderRow = $ws.Range("A43:V43")   # header-row style: col B bold+gray, rest gray, L centered/empty
$donorChangeCell = $ws.Range("L55")      # plain centered gray "no change" cell

# First row of every Datenelement group in this block (e.g. "#61" group) -
# gets the full-row gray treatment.
$headerRows = @(62, 65, 69, 73, 76, 79, 83)

# The remaining rows of each group only need the "Änderung" (L) cell reset.
$plainRows = @(63, 64, 66, 67, 68, 70, 71, 72, 74, 75, 77, 78, 80, 81, 82, 84, 85)

$donorHeaderRow.Copy()
foreach ($r in $headerRows) {
    $ws.Range("A$r`:V$r").PasteSpecial($xlPasteFormats)
}

$donorChangeCell.Copy()
foreach ($r in $headerRows) {
    $ws.Range("L$r").PasteSpecial($xlPasteFormats)
}
foreach ($r in $plainRows) {
    $ws.Range("L$r").PasteSpecial($xlPasteFormats)
}

# Clear the "ÄNDERUNG" marker text itself from column L for the whole block.
foreach ($r in ($headerRows + $plainRows)) {
    $ws.Range("L$r").Value = ""
}

$excel.CutCopyMode = $false
